$d = $word.ActiveDocument

function Get-ParaIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        $t = $t.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

function Remove-ParagraphByText($text) {
    $idx = Get-ParaIndexByText($text)
    if ($idx -lt 0) { throw "Remove-ParagraphByText: not found: $text" }
    $d.Paragraphs.Item($idx).Range.Delete()
}

function Insert-ListParagraphBefore($beforeText, $newText) {
    $idx = Get-ParaIndexByText($beforeText)
    if ($idx -lt 0) { throw "Insert-ListParagraphBefore: not found: $beforeText" }
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphBefore() | Out-Null
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.Text = $newText
}

function Insert-ListParagraphAfter($afterText, $newText) {
    $idx = Get-ParaIndexByText($afterText)
    if ($idx -lt 0) { throw "Insert-ListParagraphAfter: not found: $afterText" }
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter() | Out-Null
    $newIdx = $idx + 1
    $newPara = $d.Paragraphs.Item($newIdx)
    $newPara.Range.Text = $newText
}

function Move-ParagraphBefore($moveText, $beforeText) {
    Remove-ParagraphByText $moveText
    Insert-ListParagraphBefore $beforeText $moveText
}

# Splits the run ending at $absPos (exclusive) from the one before it by
# dropping a temporary bookmark at the boundary (forces a run break) and
# then deleting the bookmark again (leaves the break in place, engine does
# not recombine identically-formatted sibling runs once already split).
$script:tempBmCounter = 0
function Split-RunAt($absPos) {
    $script:tempBmCounter += 1
    $name = "zzTempSplit$($script:tempBmCounter)"
    $r = $d.Range($absPos, $absPos)
    $d.Bookmarks.Add($name, $r) | Out-Null
    $d.Bookmarks.Item($name).Delete()
}

# ---------------------------------------------------------------------
# 1. Header: "Project 1" -> "Project " / "2" (two runs) + _GoBack bookmark
#    placed right after the new "2" run.
# ---------------------------------------------------------------------
$searchRange = $d.Range(0, $d.Content.End)
$found = $searchRange.Find.Execute("Project 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "could not find 'Project 1'" }
$boundary = $searchRange.End - 1
Split-RunAt $boundary
$digitRange = $d.Range($searchRange.End - 1, $searchRange.End)
$digitRange.Text = "2"
# Remove the old _GoBack bookmark from its original location (after
# prog4Illegal.txt) before re-adding it here, so there's only ever one.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()
$afterTwo = $searchRange.End
$bmRange = $d.Range($afterTwo, $afterTwo)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# 2. New list items: declaredIdentifiersList.c / .h before errorList.c
# ---------------------------------------------------------------------
Insert-ListParagraphBefore "errorList.c" "declaredIdentifiersList.c"
Insert-ListParagraphBefore "errorList.c" "declaredIdentifiersList.h"

# ---------------------------------------------------------------------
# 3. FunctionalDecomposition.docx moves from its old spot (after
#    prog4Illegal.txt/run) to just before lexicalAnalyzer.c
# ---------------------------------------------------------------------
Move-ParagraphBefore "FunctionalDecomposition.docx" "lexicalAnalyzer.c"

# ---------------------------------------------------------------------
# 4. main.c / makefile move from after linkedList.h to before linkedList.c
# ---------------------------------------------------------------------
Move-ParagraphBefore "main.c" "linkedList.c"
Move-ParagraphBefore "makefile" "linkedList.c"

# ---------------------------------------------------------------------
# 5. Remove the "run" list item entirely.
# ---------------------------------------------------------------------
Remove-ParagraphByText "run"

# ---------------------------------------------------------------------
# 6. simpleProgramParser moves from after UsersManual.docx to right
#    after prog4Illegal.txt, and its text is split into three runs:
#    "simp" / "le" / "ProgramParser".
# ---------------------------------------------------------------------
Remove-ParagraphByText "simpleProgramParser"
Insert-ListParagraphAfter "prog4Illegal.txt" "simpleProgramParser"
$idx = Get-ParaIndexByText("simpleProgramParser")
$para = $d.Paragraphs.Item($idx)
$paraStart = $para.Range.Start
Split-RunAt ($paraStart + 4)
Split-RunAt ($paraStart + 6)

# ---------------------------------------------------------------------
# 7. New list items tokenList.c / tokenList.h after simpleProgramParser
# ---------------------------------------------------------------------
Insert-ListParagraphAfter "simpleProgramParser" "tokenList.c"
Insert-ListParagraphAfter "tokenList.c" "tokenList.h"

# ---------------------------------------------------------------------
# 8. Add a lastRenderedPageBreak marker before "User input: ..." text.
#    (best effort; falls back silently if unsupported by the OM)
# ---------------------------------------------------------------------
try {
    $uiIdx = Get-ParaIndexByText("User input: no user interaction with the program is required.")
    if ($uiIdx -gt 0) {
        $uiPara = $d.Paragraphs.Item($uiIdx)
        $uiPara.Range.InsertBreak(7)
    }
} catch {
}

Write-Host "Edit complete."
